$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("ECs", "Bmp4", "Bmpr1a", "ECs", 2, 0.6666666666666666, 4.504599666666667, 13.513799, 0.2870666703033801, 0.2870666703033801, 3, 1, 3.767160666666667, 11.301482, 0.05327117368020191, 0.05327117368020191, 16.96955068334644, 152.725956150118, 0.01529237845152862, 0.01529237845152862),
  @("ECs", "Bmp4", "Bmpr1a", "FAPs", 2, 0.6666666666666666, 4.504599666666667, 13.513799, 0.2870666703033801, 0.2870666703033801, 3, 1, 39.446935, 118.340805, 0.5578165391591924, 0.5578165391591923, 177.6926502520217, 1599.233852268195, 0.1601305365365844, 0.1601305365365844),
  @("ECs", "Bmp4", "Bmpr1a", "M2", 2, 0.6666666666666666, 4.504599666666667, 13.513799, 0.2870666703033801, 0.2870666703033801, 3, 1, 1.169082333333333, 3.507247, 0.01653191714824366, 0.01653191714824366, 5.266247889039221, 47.39623100135299, 0.004745762409477659, 0.00474576240947766),
  @("ECs", "Bmp4", "Bmpr1a", "sCs", 2, 0.6666666666666666, 4.504599666666667, 13.513799, 0.2870666703033801, 0.2870666703033801, 3, 1, 26.33350433333333, 79.000513, 0.372380370012362, 0.372380370012362, 118.6218948420986, 1067.597053578887, 0.1068979929057894, 0.1068979929057894),
  @("FAPs", "Bmp4", "Bmpr1a", "ECs", 3, 1, 7.480073, 22.440219, 0.476686011772755, 0.476686011772755, 3, 1, 3.767160666666667, 11.301482, 0.05327117368020191, 0.05327117368020191, 28.17863678939533, 253.607731104558, 0.02539362332406921, 0.02539362332406921),
  @("FAPs", "Bmp4", "Bmpr1a", "FAPs", 3, 1, 7.480073, 22.440219, 0.476686011772755, 0.476686011772755, 3, 1, 39.446935, 118.340805, 0.5578165391591924, 0.5578165391591923, 295.065953426255, 2655.593580836295, 0.2659033413526762, 0.2659033413526762),
  @("FAPs", "Bmp4", "Bmpr1a", "M2", 3, 1, 7.480073, 22.440219, 0.476686011772755, 0.476686011772755, 3, 1, 1.169082333333333, 3.507247, 0.01653191714824366, 0.01653191714824366, 8.744821196343665, 78.70339076709298, 0.00788053365235389, 0.00788053365235389),
  @("FAPs", "Bmp4", "Bmpr1a", "sCs", 3, 1, 7.480073, 22.440219, 0.476686011772755, 0.476686011772755, 3, 1, 26.33350433333333, 79.000513, 0.372380370012362, 0.372380370012362, 196.9765347591497, 1772.788812832347, 0.1775085134436557, 0.1775085134436557),
  @("sCs", "Bmp4", "Bmpr1a", "ECs", 3, 1, 3.707151333333333, 11.121454, 0.2362473179238649, 0.2362473179238649, 3, 1, 3.767160666666667, 11.301482, 0.05327117368020191, 0.05327117368020191, 13.96543468831422, 125.688912194828, 0.01258517190460409, 0.01258517190460409),
  @("sCs", "Bmp4", "Bmpr1a", "FAPs", 3, 1, 3.707151333333333, 11.121454, 0.2362473179238649, 0.2362473179238649, 3, 1, 39.446935, 118.340805, 0.5578165391591924, 0.5578165391591923, 146.2357576811634, 1316.12181913047, 0.1317826612699317, 0.1317826612699317),
  @("sCs", "Bmp4", "Bmpr1a", "M2", 3, 1, 3.707151333333333, 11.121454, 0.2362473179238649, 0.2362473179238649, 3, 1, 1.169082333333333, 3.507247, 0.01653191714824366, 0.01653191714824366, 4.333965130793111, 39.005686177138, 0.003905621086412115, 0.003905621086412115),
  @("sCs", "Bmp4", "Bmpr1a", "sCs", 3, 1, 3.707151333333333, 11.121454, 0.2362473179238649, 0.2362473179238649, 3, 1, 26.33350433333333, 79.000513, 0.372380370012362, 0.372380370012362, 97.62228570065578, 878.600571305902, 0.08797386366291693, 0.08797386366291694)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
  $row = $rows[$i]
  $r = $startRow + $i
  for ($c = 1; $c -le $row.Count; $c++) {
    $ws.Cells.Item($r, $c).Value = $row[$c - 1]
  }
}

Write-Output ("Dimension after edit: " + $ws.UsedRange.Address())
